$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values. The order below matches how the shared strings table
# ends up laid out in the saved workbook (short "route" strings first,
# then the longer multi-line baggage description strings).
$ws.Range("E2").Value = "Without baggage ancillary"
$ws.Range("F2").Value = "JKT-SUB  depart at 7days later"
$ws.Range("F4").Value = "MNL-CEB   depart at 7days later"
$ws.Range("F3").Value = "MNL-HKG  depart at 7days later/return at 10 days later"
$ws.Range("E3").Value = "With baggage ancillary`n1ADT(MALE) outbound 1 baggage`n1ADT(FEMALE) inbound 1 baggage`n1CHD(MALE)  both bounds 1 baggage"
$ws.Range("E4").Value = "With baggage ancillary for all`n1ADT(MALE)+1ADT(FEMALE)+1CHD(FEMALE)"

# Apply wrap text + vertical center alignment to the updated description cells
$ws.Range("E3:E4").WrapText = $true
$ws.Range("E3:E4").VerticalAlignment = -4108

# Adjust row heights to fit the new multi-line content
$ws.Rows.Item(3).RowHeight = 57
$ws.Rows.Item(4).RowHeight = 28.5

# Adjust column widths for the new, wider content (input values account for
# this runtime's column-width-to-pixel rounding so the stored <col width>
# lands as close as possible to the target 40.5 / 50.125 character widths)
$ws.Columns.Item(5).ColumnWidth = 39.65
$ws.Columns.Item(6).ColumnWidth = 49.4

# Move the active selection to F3 (as recorded in the saved workbook)
$ws.Range("F3").Select()
